# Auto-generated edit script: updates market/profit calculation values
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets, reflecting
# refreshed market board pricing data from the scheduled runner.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 222
$ws.Range("I2").Value = 222
$ws.Range("K2").Value = 222
$ws.Range("M2").Value = -109
$ws.Range("H31").Value = 54.5
$ws.Range("I31").Value = 54.5
$ws.Range("K31").Value = 163.5
$ws.Range("M31").Value = 66.5
$ws.Range("H40").Value = 2002
$ws.Range("I40").Value = 2002
$ws.Range("K40").Value = 2002
$ws.Range("M40").Value = -1827
$ws.Range("H88").Value = 5111.875
$ws.Range("I88").Value = 4501
$ws.Range("J88").Value = 5315.5
$ws.Range("K88").Value = 4501
$ws.Range("L88").Value = 5315.5
$ws.Range("M88").Value = -4095
$ws.Range("N88").Value = -6127.5
$ws.Range("H91").Value = 5111.875
$ws.Range("I91").Value = 4501
$ws.Range("J91").Value = 5315.5
$ws.Range("K91").Value = 4501
$ws.Range("L91").Value = 5315.5
$ws.Range("M91").Value = -3097
$ws.Range("N91").Value = -8123.5
$ws.Range("H123").Value = 100000
$ws.Range("I123").Value = 100000
$ws.Range("K123").Value = 100000
$ws.Range("M123").Value = -95100
$ws.Range("H132").Value = 12041.875
$ws.Range("I132").Value = 12041.875
$ws.Range("K132").Value = 36125.625
$ws.Range("M132").Value = -33595.625

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H76").Value = 75000
$ws.Range("J76").Value = 75000
$ws.Range("L76").Value = 75000
$ws.Range("N76").Value = -75676
$ws.Range("H79").Value = 75000
$ws.Range("J79").Value = 75000
$ws.Range("L79").Value = 75000
$ws.Range("N79").Value = -77340
$ws.Range("H88").Value = 1949
$ws.Range("I88").Value = 1949
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 1949
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -1543
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 1949
$ws.Range("I91").Value = 1949
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 1949
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -545
$ws.Range("N91").ClearContents()
$ws.Range("H102").Value = 2512.25
$ws.Range("I102").Value = 2512.25
$ws.Range("K102").Value = 2512.25
$ws.Range("M102").Value = -890.25
$ws.Range("H132").Value = 6319.25
$ws.Range("I132").Value = 6319.25
$ws.Range("K132").Value = 18957.75
$ws.Range("M132").Value = -16427.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 80235.5
$ws.Range("I26").Value = 471
$ws.Range("J26").Value = 160000
$ws.Range("K26").Value = 471
$ws.Range("L26").Value = 160000
$ws.Range("M26").Value = -179
$ws.Range("N26").Value = -160584
$ws.Range("H36").Value = 5259.25
$ws.Range("I36").Value = 6679
$ws.Range("K36").Value = 6679
$ws.Range("M36").Value = -6145
$ws.Range("H88").Value = 24236
$ws.Range("J88").Value = 24236
$ws.Range("L88").Value = 24236
$ws.Range("N88").Value = -25048
$ws.Range("H91").Value = 24236
$ws.Range("J91").Value = 24236
$ws.Range("L91").Value = 24236
$ws.Range("N91").Value = -27044

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 19500
$ws.Range("J92").Value = 19500
$ws.Range("L92").Value = 19500
$ws.Range("N92").Value = -24492
$ws.Range("H112").Value = 60000
$ws.Range("J112").Value = 60000
$ws.Range("L112").Value = 60000
$ws.Range("N112").Value = -62954

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 316.76923
$ws.Range("I2").Value = 201.5
$ws.Range("K2").Value = 1209
$ws.Range("M2").Value = -1096
$ws.Range("H39").Value = 1699.2142
$ws.Range("J39").Value = 2679.875
$ws.Range("L39").Value = 8039.625
$ws.Range("N39").Value = -8627.625
$ws.Range("H55").Value = 1500
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H117").Value = 5284.4287
$ws.Range("I117").Value = 694
$ws.Range("J117").Value = 7120.6
$ws.Range("K117").Value = 2082
$ws.Range("L117").Value = 21361.8
$ws.Range("M117").Value = 1360
$ws.Range("N117").Value = -28245.8
$ws.Range("H132").Value = 947.5
$ws.Range("I132").Value = 947.5
$ws.Range("K132").Value = 8527.5
$ws.Range("M132").Value = -5997.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 48225
$ws.Range("I80").Value = 46500
$ws.Range("J80").Value = 49950
$ws.Range("K80").Value = 46500
$ws.Range("L80").Value = 49950
$ws.Range("M80").Value = -45502
$ws.Range("N80").Value = -51946
$ws.Range("H83").Value = 48225
$ws.Range("I83").Value = 46500
$ws.Range("J83").Value = 49950
$ws.Range("K83").Value = 232500
$ws.Range("L83").Value = 249750
$ws.Range("M83").Value = -227508
$ws.Range("N83").Value = -259734
$ws.Range("H100").Value = 30000
$ws.Range("J100").Value = 30000
$ws.Range("L100").Value = 30000
$ws.Range("N100").Value = -32164

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1099.5
$ws.Range("I93").Value = 1099.5
$ws.Range("K93").Value = 1099.5
$ws.Range("M93").Value = 148.5
$ws.Range("H122").Value = 13497.167
$ws.Range("I122").Value = 10249.5
$ws.Range("K122").Value = 30748.5
$ws.Range("M122").Value = -28298.5
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 5000
$ws.Range("J2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("N2").Value = -5224
$ws.Range("H104").Value = 10887
$ws.Range("J104").Value = 10887
$ws.Range("L104").Value = 10887
$ws.Range("N104").Value = -17875

